# Keyboard Gloves Layout.docx - "Final changes to first working version"
#
# 1. Append a new explanatory sentence to the end of the "Overview" paragraph
#    that ends in "...position of the left-hand thumb." (adds commentary
#    about qwerty-inspired keybinds and the meaning of pinky/index finger
#    values), including two superscript ordinal suffixes ("rd", "th").
# 2. Extend the "This string of characters..." example paragraph with a
#    clause describing which character actually gets typed.
# 3/4. Mark the two table cells where content now lands on a fresh page
#    (after the text above pushes things down) with <w:lastRenderedPageBreak/>.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Extend the paragraph ending "...position of the left-hand thumb."
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(4)
$run1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)

$plainRuns1 = @(
    " ",
    "Also note that the keybinds ",
    "took inspiration",
    " ",
    "from",
    " on a traditional qwerty keyboard layout",
    " to help the user",
    " have an easy transition from a ",
    "qwerty keyboard to using the Keyboard Gloves.",
    " ",
    "A",
    " ",
    "value of 1 for the",
    " ",
    "left hand",
    " pinky represents the left side of the keyboard and a value of 1 for the ",
    "left",
    " hand index finger represents the right side of the keyboard",
    ", with "
)
foreach ($t in $plainRuns1) {
    $run1.InsertAfter($t)
}

# "the 3" + superscript "rd"
$run1.InsertAfter("the 3")
$supStart1 = $run1.End
$run1.InsertAfter("rd")
$d.Range($supStart1, $run1.End).Font.Superscript = $true

# " and 4" + superscript "th"
$run1.InsertAfter(" and 4")
$supStart2 = $run1.End
$run1.InsertAfter("th")
$d.Range($supStart2, $run1.End).Font.Superscript = $true

# closing text
$run1.InsertAfter(" fingers specifying the row of keys that the user wants to type from.")

# ---------------------------------------------------------------------------
# 2) Extend the "This string of characters..." example paragraph
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This string of characters will result from bending the left-hand index finger and the right-hand ring finger, with all other fingers extended.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This string of characters will result from bending the left-hand index finger and the right-hand ring finger, with all other fingers extended",
    2) | Out-Null

$p2 = $d.Paragraphs.Item(7)
$run2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$run2.InsertAfter([string][char]0x0020 + "and will type the character " + [string][char]0x201C)
$run2.InsertAfter("o" + [string][char]0x201D + " being typed.")

# ---------------------------------------------------------------------------
# 3/4) Mark the two table cells that now start a fresh rendered page with
#      <w:lastRenderedPageBreak/> as the first child of their first run.
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)

function Add-LastRenderedPageBreak($cellRange, $innerXml) {
    $full = $d.Range($cellRange.Start, $cellRange.End - 1)
    $pkg = '<?xml version="1.0"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($pkg) | Out-Null
}

# Row with first-column value "11100" (first occurrence in the table).
$cellA = $tbl.Cell(25, 1).Range
Add-LastRenderedPageBreak $cellA '<w:p><w:r><w:lastRenderedPageBreak/><w:t>11100</w:t></w:r></w:p>'

# Row with first-column value "0011" + "1" (first occurrence in the table).
$cellB = $tbl.Cell(70, 1).Range
Add-LastRenderedPageBreak $cellB '<w:p><w:r><w:lastRenderedPageBreak/><w:t>0011</w:t></w:r><w:r><w:t>1</w:t></w:r></w:p>'

Write-Output "Edit complete"
